$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'47.330.15"
$ws.Range("E2").Value = "'  +0.63%  "
$ws.Range("D3").Value = "'2.489.44"
$ws.Range("E3").Value = "'  -0.22%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'321.18"
$ws.Range("E5").Value = "'  -0.35%  "
$ws.Range("D6").Value = "'108.55"
$ws.Range("E6").Value = "'  +2.86%  "
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = "'  -0.93%  "
$ws.Range("D10").Value = "'39.41"
$ws.Range("E10").Value = "'  +5.06%  "
$ws.Range("E11").Value = "'  -0.47%  "
$ws.Range("E12").Value = "'  +0.39%  "
$ws.Range("D13").Value = "'18.34"
$ws.Range("E13").Value = "'  -0.11%  "
$ws.Range("E14").Value = "'  -1.20%  "
$ws.Range("D15").Value = "'2.881.47"
$ws.Range("E15").Value = "'  -0.14%  "
$ws.Range("D16").Value = "'2.481.16"
$ws.Range("E16").Value = "'  -0.44%  "
$ws.Range("D17").Value = "'0.845"
$ws.Range("E17").Value = "'  -0.04%  "
$ws.Range("D18").Value = "'47.217.26"
$ws.Range("E18").Value = "'  +0.50%  "
$ws.Range("D19").Value = "'13.14"
$ws.Range("E19").Value = "'  +3.98%  "
$ws.Range("E20").Value = "'  +1.13%  "
$ws.Range("D21").Value = "'0.0₃0936"
$ws.Range("E21").Value = "'  +0.25%  "
$ws.Range("D22").Value = "'2.65"
$ws.Range("E22").Value = "'  +12.46%  "
$ws.Range("E23").Value = "'  -0.77%  "
$ws.Range("E24").Value = "'  -2.60%  "
$ws.Range("E25").Value = "'  +0.55%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "'  -0.09%  "
$ws.Range("D27").Value = "'25.72"
$ws.Range("E27").Value = "'  -1.58%  "
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "'  +3.12%  "
$ws.Range("D29").Value = "'9.99"
$ws.Range("E29").Value = "'  -1.74%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.138"
$ws.Range("E30").Value = "'  +2.16%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'34.68"
$ws.Range("E31").Value = "'  -1.62%  "
$ws.Range("E32").Value = "'  +0.30%  "
$ws.Range("D33").Value = "'20.57"
$ws.Range("E33").Value = "'  +4.00%  "
$ws.Range("E34").Value = "'  -0.53%  "
$ws.Range("E35").Value = "'  +0.13%  "
$ws.Range("D36").Value = "'1.01"
$ws.Range("E36").Value = "'  +0.23%  "
$ws.Range("D37").Value = "'4.74"
$ws.Range("E37").Value = "'  +2.76%  "
$ws.Range("D39").Value = "'2.93"
$ws.Range("E39").Value = "'  -1.74%  "
$ws.Range("D40").Value = "'22.96"
$ws.Range("E40").Value = "'  +5.55%  "
$ws.Range("E41").Value = "'  -0.01%  "
$ws.Range("E42").Value = "'  +0.48%  "
$ws.Range("D43").Value = "'116.79"
$ws.Range("E43").Value = "'  -4.35%  "
$ws.Range("E44").Value = "'  +0.59%  "
$ws.Range("D45").Value = "'1.997.70"
$ws.Range("E45").Value = "'  +2.27%  "
$ws.Range("E46").Value = "'  +2.16%  "
$ws.Range("E47").Value = "'  -5.07%  "
$ws.Range("D48").Value = "'9.18"
$ws.Range("E48").Value = "'  +0.20%  "
$ws.Range("E49").Value = "'  -0.58%  "
$ws.Range("D50").Value = "'5.09"
$ws.Range("E50").Value = "'  -5.15%  "
$ws.Range("D51").Value = "'56.45"
$ws.Range("E51").Value = "'  +3.31%  "

Write-Output "Applied cryptos update"
